$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Feria Lagunitas de Puerto
# Montt" / Espinaca. It belongs right before the current row 39, so insert a
# fresh row there; this pushes the existing rows 39:101 down to 40:102
# (the old last row, 101, ends up duplicated/shifted to the new last row, 102).
$ws.Rows("39:39").Insert()

# Fill the newly inserted row 39 using the same constant template shared by
# every other data row in this sheet, together with the new observation's
# date, volume and price figures.
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 45272
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = 100112012
$ws.Cells.Item(39, 7).Value = "Espinaca"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 40
$ws.Cells.Item(39, 11).Value = 15000
$ws.Cells.Item(39, 12).Value = 15000
$ws.Cells.Item(39, 13).Value = 15000
$ws.Cells.Item(39, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(39, 15).Value = "Región Metropolitana"
$ws.Cells.Item(39, 16).Value = 1500
$ws.Cells.Item(39, 17).Value = 10
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# Keep the date cell's number format consistent with the other date cells.
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(40, 4).NumberFormat
